# Auto-generated edit script: updates market-price derived cells (H:N)
# across multiple worksheets, matching the scheduled runner's refreshed values.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1975.9697
$ws.Range("I19").Value = 3536.75
$ws.Range("J19").Value = 507
$ws.Range("K19").Value = 3536.75
$ws.Range("L19").Value = 507
$ws.Range("M19").Value = -3361.75
$ws.Range("N19").Value = -857
$ws.Range("H118").Value = 1899
$ws.Range("J118").Value = 2817.25
$ws.Range("L118").Value = 8451.75
$ws.Range("N118").Value = -11765.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N10").ClearContents()
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H74").Value = 1646.762
$ws.Range("I74").Value = 1097.74
$ws.Range("J74").Value = 3758.3845
$ws.Range("K74").Value = 1097.74
$ws.Range("L74").Value = 3758.3845
$ws.Range("M74").Value = -223.74
$ws.Range("N74").Value = -5506.3845
$ws.Range("H77").Value = 1646.762
$ws.Range("I77").Value = 1097.74
$ws.Range("J77").Value = 3758.3845
$ws.Range("K77").Value = 5488.7
$ws.Range("L77").Value = 18791.9225
$ws.Range("M77").Value = -1120.7
$ws.Range("N77").Value = -27527.9225

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2602.4707
$ws.Range("I20").Value = 2527.5
$ws.Range("J20").Value = 2643.3635
$ws.Range("K20").Value = 2527.5
$ws.Range("L20").Value = 2643.3635
$ws.Range("M20").Value = -2280.5
$ws.Range("N20").Value = -3137.3635
$ws.Range("H64").Value = 408.3125
$ws.Range("I64").Value = 155
$ws.Range("J64").Value = 523.4545000000001
$ws.Range("K64").Value = 155
$ws.Range("L64").Value = 523.4545000000001
$ws.Range("M64").Value = 70
$ws.Range("N64").Value = -973.4545000000001
$ws.Range("H67").Value = 408.3125
$ws.Range("I67").Value = 155
$ws.Range("J67").Value = 523.4545000000001
$ws.Range("K67").Value = 155
$ws.Range("L67").Value = 523.4545000000001
$ws.Range("M67").Value = 625
$ws.Range("N67").Value = -2083.4545
$ws.Range("H99").Value = 5328175
$ws.Range("I99").Value = 1834721.2
$ws.Range("J99").Value = 20000680
$ws.Range("K99").Value = 1834721.2
$ws.Range("L99").Value = 20000680
$ws.Range("M99").Value = -1833223.2
$ws.Range("N99").Value = -20003676

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N16").Value = -1807.3334
$ws.Range("H16").Value = 1971.1
$ws.Range("I16").Value = 2287.2856
$ws.Range("J16").Value = 1233.3334
$ws.Range("K16").Value = 2287.2856
$ws.Range("L16").Value = 1233.3334
$ws.Range("M16").Value = -2000.2856
$ws.Range("H22").Value = 791.5
$ws.Range("I22").Value = 482.3
$ws.Range("J22").Value = 1100.7
$ws.Range("K22").Value = 482.3
$ws.Range("L22").Value = 1100.7
$ws.Range("M22").Value = -132.3
$ws.Range("N22").Value = -1800.7
$ws.Range("H31").Value = 2430.7441
$ws.Range("I31").Value = 1701.7291
$ws.Range("J31").Value = 3351.6052
$ws.Range("K31").Value = 1701.7291
$ws.Range("L31").Value = 3351.6052
$ws.Range("M31").Value = -1406.7291
$ws.Range("N31").Value = -3941.6052
$ws.Range("H34").Value = 2430.7441
$ws.Range("I34").Value = 1701.7291
$ws.Range("J34").Value = 3351.6052
$ws.Range("K34").Value = 1701.7291
$ws.Range("L34").Value = 3351.6052
$ws.Range("M34").Value = -1499.7291
$ws.Range("N34").Value = -3755.6052
$ws.Range("H112").Value = 29910
$ws.Range("J112").Value = 29910
$ws.Range("L112").Value = 29910
$ws.Range("N112").Value = -32864
$ws.Range("N113").Value = -5573.3334
$ws.Range("H113").Value = 1971.1
$ws.Range("I113").Value = 2287.2856
$ws.Range("J113").Value = 1233.3334
$ws.Range("K113").Value = 2287.2856
$ws.Range("L113").Value = 1233.3334
$ws.Range("M113").Value = -117.2856000000002
$ws.Range("H132").Value = 2307.0466
$ws.Range("I132").Value = 1072.7727
$ws.Range("J132").Value = 3600.0952
$ws.Range("K132").Value = 3218.3181
$ws.Range("L132").Value = 10800.2856
$ws.Range("M132").Value = -688.3181
$ws.Range("N132").Value = -15860.2856
$ws.Range("H134").Value = 1765.5834
$ws.Range("I134").Value = 1029.5264
$ws.Range("J134").Value = 2588.2354
$ws.Range("K134").Value = 3088.5792
$ws.Range("L134").Value = 7764.706200000001
$ws.Range("M134").Value = -553.5792000000001
$ws.Range("N134").Value = -12834.7062

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 521.8889
$ws.Range("I68").Value = 487.7143
$ws.Range("J68").Value = 543.63635
$ws.Range("K68").Value = 1463.1429
$ws.Range("L68").Value = 1630.90905
$ws.Range("M68").Value = -652.1428999999998
$ws.Range("N68").Value = -3252.90905
$ws.Range("H71").Value = 521.8889
$ws.Range("I71").Value = 487.7143
$ws.Range("J71").Value = 543.63635
$ws.Range("K71").Value = 4389.428699999999
$ws.Range("L71").Value = 4892.72715
$ws.Range("M71").Value = -333.4286999999995
$ws.Range("N71").Value = -13004.72715
$ws.Range("H86").Value = 956.6667
$ws.Range("I86").Value = 110
$ws.Range("J86").Value = 2650
$ws.Range("K86").Value = 330
$ws.Range("L86").Value = 7950
$ws.Range("M86").Value = 856
$ws.Range("N86").Value = -10322
$ws.Range("H89").Value = 956.6667
$ws.Range("I89").Value = 110
$ws.Range("J89").Value = 2650
$ws.Range("K89").Value = 990
$ws.Range("L89").Value = 23850
$ws.Range("M89").Value = 4938
$ws.Range("N89").Value = -35706

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 12503751
$ws.Range("I18").Value = 12503751
$ws.Range("K18").Value = 12503751
$ws.Range("M18").Value = -12503458
$ws.Range("H102").Value = 4471.174
$ws.Range("I102").Value = 5073.5557
$ws.Range("K102").Value = 5073.5557
$ws.Range("M102").Value = -3451.5557

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 461
$ws.Range("I22").Value = 398.6
$ws.Range("K22").Value = 398.6
$ws.Range("M22").Value = -103.6
$ws.Range("H27").Value = 461
$ws.Range("I27").Value = 398.6
$ws.Range("K27").Value = 398.6
$ws.Range("M27").Value = -291.6

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M14").ClearContents()
$ws.Range("H14").Value = 34336332
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 34336332
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 34336332
$ws.Range("N14").Value = -34336668
